# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 0
    4  = 4
    5  = 5
    6  = 3
    7  = 0
    8  = 1
    9  = 1
    10 = 3
    11 = 5
    12 = 1
    13 = 1
    14 = 7
    15 = 1
    16 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
